$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix tiny precision errors in existing values
$ws.Cells.Item(1351, 2).Value2 = 0.9456975608411728
$ws.Cells.Item(1352, 2).Value2 = 0.9456215060157313
$ws.Cells.Item(1353, 2).Value2 = 0.9457021287033124
$ws.Cells.Item(1355, 2).Value2 = 0.9520866630861704
$ws.Cells.Item(1359, 2).Value2 = 0.9444399514788177

# Append new rows of currency price data
$dates = @(
    "2023-10-09", "2023-10-10", "2023-10-11", "2023-10-12", "2023-10-13", "2023-10-14", "2023-10-15", "2023-10-16", "2023-10-17", "2023-10-18", "2023-10-19", "2023-10-20", "2023-10-21", "2023-10-22", "2023-10-23", "2023-10-24", "2023-10-25", "2023-10-26", "2023-10-27", "2023-10-28", "2023-10-29", "2023-10-30", "2023-10-31", "2023-11-01", "2023-11-02", "2023-11-03", "2023-11-04", "2023-11-05", "2023-11-06", "2023-11-07", "2023-11-08", "2023-11-09", "2023-11-10", "2023-11-11", "2023-11-12", "2023-11-13", "2023-11-14", "2023-11-15", "2023-11-16", "2023-11-17", "2023-11-18", "2023-11-19", "2023-11-20", "2023-11-21", "2023-11-22", "2023-11-23", "2023-11-24", "2023-11-25", "2023-11-26", "2023-11-27", "2023-11-28", "2023-11-29", "2023-11-30", "2023-12-01", "2023-12-02", "2023-12-03", "2023-12-04", "2023-12-05", "2023-12-06", "2023-12-07", "2023-12-08", "2023-12-09", "2023-12-10", "2023-12-11", "2023-12-12", "2023-12-13", "2023-12-14", "2023-12-15", "2023-12-16", "2023-12-17", "2023-12-18", "2023-12-19", "2023-12-20", "2023-12-21", "2023-12-22", "2023-12-23", "2023-12-24", "2023-12-25", "2023-12-26", "2023-12-27", "2023-12-28", "2023-12-29", "2023-12-30", "2023-12-31", "2024-01-01", "2024-01-02", "2024-01-03", "2024-01-04", "2024-01-05", "2024-01-06", "2024-01-07", "2024-01-08", "2024-01-09", "2024-01-10", "2024-01-11", "2024-01-12", "2024-01-13", "2024-01-14", "2024-01-15", "2024-01-16", "2024-01-17", "2024-01-18", "2024-01-19", "2024-01-20", "2024-01-21", "2024-01-22", "2024-01-23", "2024-01-24", "2024-01-25", "2024-01-26", "2024-01-27", "2024-01-28", "2024-01-29", "2024-01-30", "2024-01-31", "2024-02-01", "2024-02-02", "2024-02-03", "2024-02-04", "2024-02-05", "2024-02-06", "2024-02-07", "2024-02-08", "2024-02-09", "2024-02-10", "2024-02-11", "2024-02-12", "2024-02-13", "2024-02-14", "2024-02-15", "2024-02-16", "2024-02-17", "2024-02-18", "2024-02-19", "2024-02-20", "2024-02-21", "2024-02-22", "2024-02-23", "2024-02-24"
)
$prices = @(
    0.9457616408940103, 0.9458719939969993, 0.9418504513989946, 0.9371815162452983, 0.9486430923517377, 0.9479669289916347, 0.9517374997653508, 0.9511400879149025, 0.9466280156956492, 0.9465662618061806, 0.9484659787173039, 0.9451491541775148, 0.9447987159158298, 0.9447355074834111, 0.9442536347444758, 0.9368670245318715, 0.9426736264548318, 0.9460354420237632, 0.9484597745663718, 0.9465690694959129, 0.9458378644196922, 0.9456712083248112, 0.9428466669548856, 0.9461954440959768, 0.9464282241132318, 0.9422328909815197, 0.9312958220851334, 0.9347413209698966, 0.93599725342184, 0.9329003033454832, 0.9398552860879338, 0.9346188085283772, 0.9382198088392683, 0.9351108950851816, 0.9394146981881806, 0.9351248221175344, 0.935437059715657, 0.9198884362648033, 0.9227007321539276, 0.9211450709065637, 0.9172760212037482, 0.9163590501137914, 0.9178315205840671, 0.9152102573227366, 0.91780170953887, 0.920798952059709, 0.9190593618577992, 0.9165332460469506, 0.9136846214640336, 0.9172627317194424, 0.9129141754175456, 0.9107151448786911, 0.9118489704293847, 0.9183369929718712, 0.9245729135096348, 0.9252191436011231, 0.9223545575311097, 0.9255011484545527, 0.9270924392555342, 0.9279949122675909, 0.9330513399330395, 0.9310799751007045, 0.928719676334498, 0.9305730693638169, 0.9297222741172922, 0.9274348028709091, 0.9253544878069694, 0.9193573211760886, 0.9212535261235121, 0.9217396584468978, 0.9196428927335351, 0.9172514027250903, 0.9167044391877737, 0.9133102709048582, 0.9195751568258463, 0.9174826020863721, 0.9145402601139501, 0.9211659525050861, 0.9130638753410724, 0.9071328404268897, 0.9150287467751845, 0.9067552655660333, 0.9146399696325263, 0.9119386801207846, 0.9107184149179123, 0.9129065976846708, 0.9159844082939458, 0.917466740137594, 0.918839045971852, 0.9242187922828472, 0.9206507884362776, 0.9217330541492468, 0.9143213678229255, 0.9156338051625903, 0.9139927580372474, 0.9112572388837912, 0.9119137373957775, 0.9136861486895842, 0.9121919062818314, 0.9179701273478433, 0.9178831711712865, 0.9169759046663901, 0.9173549116085207, 0.9167272155977474, 0.9170967013011817, 0.9192971552632508, 0.915402982469925, 0.921190806539734, 0.9186637115475811, 0.9215752099770762, 0.9220461565638147, 0.9226640195107103, 0.9231873702579477, 0.922320793323699, 0.9238868222587213, 0.9237851656620868, 0.9229203745253124, 0.9258987421882034, 0.9258508526980923, 0.9275778065221452, 0.9306844907061551, 0.9319862797869902, 0.928160669022553, 0.9304116073029959, 0.9270274560368433, 0.9270975523262067, 0.9284688390290137, 0.9284746004560016, 0.9343293487007063, 0.9318769402412725, 0.9299918562438134, 0.9307126662581066, 0.9279043383810303, 0.9270099574161197, 0.9302225599850594, 0.9262042734330065, 0.9254005408731009, 0.9246029106702607, 0.9268059536371416
)

$startRow = 1361
for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $startRow + $i
    $cellA = $ws.Cells.Item($row, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value2 = $dates[$i]
    $cellA.Style = "Normal"
    $ws.Cells.Item($row, 2).Value2 = $prices[$i]
}
